# Update cryptocurrency price/volume figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.669.72"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").Value = "1.598.41"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").Value = "'211.54"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("E6").Value = "  +0.78%  "

$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").Value = "'0.0619"
$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("E9").Value = "  +0.73%  "

$ws.Range("D10").Value = "'19.58"
$ws.Range("E10").Value = "  -0.50%  "

$ws.Range("E11").Value = "  +0.33%  "

$ws.Range("D12").Value = "1.822.14"

$ws.Range("D13").Value = "1.591.42"
$ws.Range("E13").Value = "  +1.31%  "

$ws.Range("E15").Value = "  +0.51%  "

$ws.Range("D16").Value = "'65.20"
$ws.Range("E16").Value = "  +0.43%  "

$ws.Range("D17").Value = "26.658.00"
$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("D18").Value = "0.0₃0737"
$ws.Range("E18").Value = "  +1.39%  "

$ws.Range("D19").Value = "'209.60"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D21").Value = "'7.03"
$ws.Range("E21").Value = "  +3.86%  "

$ws.Range("D22").Value = "'4.29"
$ws.Range("E22").Value = "  +0.71%  "

$ws.Range("E23").Value = "  +1.68%  "

$ws.Range("E24").Value = "  +0.79%  "

$ws.Range("D25").Value = "'144.32"
$ws.Range("E25").Value = "  -1.50%  "

$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("D27").Value = "'7.13"
$ws.Range("E27").Value = "  -0.50%  "

$ws.Range("E28").Value = "  -0.32%  "

$ws.Range("D29").Value = "'15.29"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("E30").Value = "  +2.35%  "

$ws.Range("E31").Value = "  +0.08%  "

$ws.Range("E32").Value = "  +0.81%  "

$ws.Range("E33").Value = "  +1.64%  "

$ws.Range("D34").Value = "1.289.28"
$ws.Range("E34").Value = "  -0.75%  "

$ws.Range("E35").Value = "  -6.76%  "

$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("E37").Value = "  +0.70%  "

$ws.Range("E38").Value = "  -0.58%  "

$ws.Range("E39").Value = "  -0.89%  "

$ws.Range("D40").Value = "'1.05"
$ws.Range("E40").Value = "  +19.51%  "

$ws.Range("E41").Value = "  +2.49%  "

$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("D43").Value = "'0.784"
$ws.Range("E43").Value = "  -0.71%  "

$ws.Range("D44").Value = "'63.61"
$ws.Range("E44").Value = "  -0.36%  "

$ws.Range("D45").Value = "1.735.81"
$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").Value = "'90.77"
$ws.Range("E46").Value = "  +0.72%  "

$ws.Range("E47").Value = "  -3.31%  "

$ws.Range("E48").Value = "  -0.06%  "

$ws.Range("E49").Value = "  +1.66%  "

$ws.Range("D50").Value = "'0.0508"
$ws.Range("E50").Value = "  +0.84%  "

$ws.Range("E51").Value = "  +0.32%  "
